# Update bulk upload template to use passport codes instead of player names
# Sheet 1: "比赛数据" (Match Data) - replace Chinese player names with passport codes,
#          and replace full-width gender text with single-letter English codes.
# Sheet 2: "使用说明" (Instructions) - rewrite instruction text to reference passport codes.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 比赛数据 ----
$ws1 = $wb.Worksheets.Item("比赛数据")

$ws1.Range("A2").Value = "HVGN0BW0"
$ws1.Range("C2").Value = "KGLE38K4"

$ws1.Range("A3").Value = "MJST45X9"
$ws1.Range("B3").Value = "SWQR78Z2"
$ws1.Range("C3").Value = "TBPL91M5"
$ws1.Range("D3").Value = "LCKM33Y8"

$ws1.Range("A4").Value = "DLVW67N4"
$ws1.Range("C4").Value = "EDRX29H6"
$ws1.Range("H4").Value = "M"

$ws1.Range("A5").Value = "AWJF82P1"
$ws1.Range("B5").Value = "JLMN56Q3"
$ws1.Range("C5").Value = "RZQW74T7"
$ws1.Range("D5").Value = "MKHY93V0"

$ws1.Range("A6").Value = "CTBR48K5"
$ws1.Range("C6").Value = "ARLZ61F9"
$ws1.Range("H6").Value = "F"

# ---- Sheet 2: 使用说明 ----
$ws2 = $wb.Worksheets.Item("使用说明")

$ws2.Range("A3").Value = "格式指南："
$ws2.Range("A4").Value = "• 使用选手护照代码（例如：HVGN0BW0, KGLE38K4）"
$ws2.Range("A5").Value = "• 单打比赛请将第一队选手二和第二队选手二留空"
$ws2.Range("A6").Value = "• 使用如下分数格式：11, 7, 15, 13（游戏比分）"
$ws2.Range("A7").Value = "• 日期格式：YYYY-MM-DD（例如：2025-01-15）"
$ws2.Range("A8").Value = "• 性别覆盖：男或女（可选，用于跨性别比赛）"

$ws2.Range("A10").Value = "示例："
$ws2.Range("A11").Value = "单打：HVGN0BW0 对 KGLE38K4，比分 11-7"
$ws2.Range("A12").Value = "双打：MJST45X9/SWQR78Z2 对 TBPL91M5/LCKM33Y8，比分 11-9"

# Row 13 becomes a blank line in the instructions (kept as an explicit empty
# text cell, matching its siblings A2/A9, rather than a fully cleared cell).
$ws2.Range("A13").Value = "'"
$ws2.Range("A13").Style = "Normal"

$ws2.Range("A14").Value = "验证将检查："
$ws2.Range("A15").Value = "• 所有护照代码在系统中存在"
$ws2.Range("A16").Value = "• 有效的分数格式"
$ws2.Range("A17").Value = "• 无重复比赛"
$ws2.Range("A18").Value = "• 正确的日期格式"
